$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3257.1794
$ws.Range("I40").Value = 3603.9167
$ws.Range("K40").Value = 3603.9167
$ws.Range("M40").Value = -3428.9167
$ws.Range("H58").Value = 4563.385
$ws.Range("J58").Value = 5245.8184
$ws.Range("L58").Value = 15737.4552
$ws.Range("N58").Value = -16037.4552
$ws.Range("H76").Value = 11555.556
$ws.Range("J76").Value = 7985.7144
$ws.Range("L76").Value = 7985.7144
$ws.Range("N76").Value = -8615.714400000001
$ws.Range("H79").Value = 11555.556
$ws.Range("J79").Value = 7985.7144
$ws.Range("L79").Value = 7985.7144
$ws.Range("N79").Value = -10169.7144
$ws.Range("H86").Value = 3942.0557
$ws.Range("I86").Value = 2844.75
$ws.Range("J86").Value = 4819.9
$ws.Range("K86").Value = 2844.75
$ws.Range("L86").Value = 4819.9
$ws.Range("M86").Value = -1721.75
$ws.Range("N86").Value = -7065.9
$ws.Range("H89").Value = 3942.0557
$ws.Range("I89").Value = 2844.75
$ws.Range("J89").Value = 4819.9
$ws.Range("K89").Value = 14223.75
$ws.Range("L89").Value = 24099.5
$ws.Range("M89").Value = -8607.75
$ws.Range("N89").Value = -35331.5
$ws.Range("H113").Value = 5287.2856
$ws.Range("I113").Value = 4613.909
$ws.Range("K113").Value = 4613.909
$ws.Range("M113").Value = -1359.909
$ws.Range("H135").Value = 2739.1428
$ws.Range("I135").Value = 2895.6667
$ws.Range("J135").Value = 1800
$ws.Range("K135").Value = 26061.0003
$ws.Range("L135").Value = 16200
$ws.Range("M135").Value = -23526.0003
$ws.Range("N135").Value = -21270

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9131.33
$ws.Range("I32").Value = 3033.84
$ws.Range("J32").Value = 27423.8
$ws.Range("K32").Value = 3033.84
$ws.Range("L32").Value = 27423.8
$ws.Range("M32").Value = -2746.84
$ws.Range("N32").Value = -27997.8
$ws.Range("H97").Value = 3515.889
$ws.Range("I97").Value = 3300.4666
$ws.Range("K97").Value = 3300.4666
$ws.Range("M97").Value = -2804.4666
$ws.Range("H102").Value = 17744.111
$ws.Range("I102").Value = 17744.111
$ws.Range("K102").Value = 17744.111
$ws.Range("M102").Value = -16122.111
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5110.0625
$ws.Range("I20").Value = 4832.923
$ws.Range("J20").Value = 6311
$ws.Range("K20").Value = 4832.923
$ws.Range("L20").Value = 6311
$ws.Range("M20").Value = -4585.923
$ws.Range("N20").Value = -6805
$ws.Range("H94").Value = 2295.0435
$ws.Range("I94").Value = 2859.2
$ws.Range("J94").Value = 2138.3333
$ws.Range("K94").Value = 2859.2
$ws.Range("L94").Value = 2138.3333
$ws.Range("M94").Value = -2408.2
$ws.Range("N94").Value = -3040.3333
$ws.Range("H99").Value = 5541.3228
$ws.Range("I99").Value = 6608.8945
$ws.Range("K99").Value = 6608.8945
$ws.Range("M99").Value = -5110.8945
$ws.Range("H134").Value = 3108423.2
$ws.Range("I134").Value = 1721.1613
$ws.Range("J134").Value = 9528941
$ws.Range("K134").Value = 5163.4839
$ws.Range("L134").Value = 28586823
$ws.Range("M134").Value = -2628.4839
$ws.Range("N134").Value = -28591893

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2992.5
$ws.Range("J99").Value = 3135.25
$ws.Range("L99").Value = 3135.25
$ws.Range("N99").Value = -6131.25
$ws.Range("H122").Value = 2461
$ws.Range("I122").Value = 2502.2
$ws.Range("K122").Value = 7506.599999999999
$ws.Range("M122").Value = -5056.599999999999
$ws.Range("H126").Value = 2992.5
$ws.Range("J126").Value = 3135.25
$ws.Range("L126").Value = 9405.75
$ws.Range("N126").Value = -14345.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 19800
$ws.Range("I87").Value = 13000
$ws.Range("K87").Value = 39000
$ws.Range("M87").Value = -37752
$ws.Range("H90").Value = 19800
$ws.Range("I90").Value = 13000
$ws.Range("K90").Value = 117000
$ws.Range("M90").Value = -110760
$ws.Range("H92").Value = 348
$ws.Range("I92").Value = 299
$ws.Range("K92").Value = 897
$ws.Range("M92").Value = 351
$ws.Range("H94").Value = 2947.913
$ws.Range("I94").Value = 828.8570999999999
$ws.Range("K94").Value = 2486.5713
$ws.Range("M94").Value = -1810.5713
$ws.Range("H139").Value = 4763.5835
$ws.Range("I139").Value = 4527.8335
$ws.Range("K139").Value = 13583.5005
$ws.Range("M139").Value = -8443.500499999998
$ws.Range("H141").Value = 607996
$ws.Range("I141").Value = 607996
$ws.Range("K141").Value = 1823988
$ws.Range("M141").Value = -1818808

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5764.923
$ws.Range("I70").Value = 5408
$ws.Range("K70").Value = 5408
$ws.Range("M70").Value = -5138
$ws.Range("H73").Value = 5764.923
$ws.Range("I73").Value = 5408
$ws.Range("K73").Value = 5408
$ws.Range("M73").Value = -4472
$ws.Range("H97").Value = 1555.6786
$ws.Range("I97").Value = 1133.2632
$ws.Range("K97").Value = 1133.2632
$ws.Range("M97").Value = -637.2632000000001
$ws.Range("H113").Value = 4006.92
$ws.Range("I113").Value = 3384.1875
$ws.Range("J113").Value = 5114
$ws.Range("K113").Value = 3384.1875
$ws.Range("L113").Value = 5114
$ws.Range("M113").Value = -1214.1875
$ws.Range("N113").Value = -9454
$ws.Range("H126").Value = 4625
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H46").Value = 3917
$ws.Range("I46").Value = 2996.25
$ws.Range("J46").Value = 4735.4443
$ws.Range("K46").Value = 2996.25
$ws.Range("L46").Value = 4735.4443
$ws.Range("M46").Value = -2808.25
$ws.Range("N46").Value = -5111.4443
$ws.Range("H93").Value = 83335560
$ws.Range("I93").Value = 111113100
$ws.Range("J93").Value = 2929.6667
$ws.Range("K93").Value = 111113100
$ws.Range("L93").Value = 2929.6667
$ws.Range("M93").Value = -111111852
$ws.Range("N93").Value = -5425.6667
$ws.Range("H100").Value = 1896.1428
$ws.Range("I100").Value = 1445.875
$ws.Range("J100").Value = 2496.5
$ws.Range("K100").Value = 1445.875
$ws.Range("L100").Value = 2496.5
$ws.Range("M100").Value = -904.875
$ws.Range("N100").Value = -3578.5
$ws.Range("H132").Value = 4985.067
$ws.Range("I132").Value = 5862
$ws.Range("J132").Value = 3982.8572
$ws.Range("K132").Value = 17586
$ws.Range("L132").Value = 11948.5716
$ws.Range("M132").Value = -15056
$ws.Range("N132").Value = -17008.5716
$ws.Range("H136").Value = 51518.242
$ws.Range("J136").Value = 142234.22
$ws.Range("L136").Value = 426702.66
$ws.Range("N136").Value = -431802.66

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12029.8125
$ws.Range("I81").Value = 5416.1665
$ws.Range("J81").Value = 15998
$ws.Range("K81").Value = 10832.333
$ws.Range("L81").Value = 31996
$ws.Range("M81").Value = -9771.333000000001
$ws.Range("N81").Value = -34118
$ws.Range("H84").Value = 12029.8125
$ws.Range("I84").Value = 5416.1665
$ws.Range("J84").Value = 15998
$ws.Range("K84").Value = 54161.665
$ws.Range("L84").Value = 159980
$ws.Range("M84").Value = -48857.665
$ws.Range("N84").Value = -170588
$ws.Range("H100").Value = 1136
$ws.Range("I100").Value = 1056.4286
$ws.Range("K100").Value = 2112.8572
$ws.Range("M100").Value = -1571.8572
$ws.Range("H107").Value = 16130235
$ws.Range("I107").Value = 21740270
$ws.Range("J107").Value = 1382.375
$ws.Range("K107").Value = 65220810
$ws.Range("L107").Value = 4147.125
$ws.Range("M107").Value = -65218890
$ws.Range("N107").Value = -7987.125
$ws.Range("H113").Value = 259.66666
$ws.Range("I113").Value = 259.66666
$ws.Range("K113").Value = 778.9999799999999
$ws.Range("M113").Value = 1391.00002
$ws.Range("H122").Value = 3109.4443
$ws.Range("I122").Value = 3109.4443
$ws.Range("K122").Value = 9328.332900000001
$ws.Range("M122").Value = -6878.332900000001
$ws.Range("H126").Value = 1485.9231
$ws.Range("I126").Value = 1479.1765
$ws.Range("K126").Value = 4437.529500000001
$ws.Range("M126").Value = -1967.529500000001
$ws.Range("H132").Value = 5214802
$ws.Range("I132").Value = 6358.4
$ws.Range("J132").Value = 13895541
$ws.Range("K132").Value = 19075.2
$ws.Range("L132").Value = 41686623
$ws.Range("M132").Value = -16545.2
$ws.Range("N132").Value = -41691683
